# Accelerate process runner and accelerate purchasers.
#
# Sheet "Лист1" layout (row -> meaning):
#   Row 2: "Закупка материалов" (purchasing materials)  -> D2 time in seconds
#   Row 3: "Поставка материалов" (delivery of materials) -> D3 time in seconds
#   Row 4: "Резка штрипса" (strip cutting)               -> D4 time in seconds
#   Row 17: "Между запусками в 1 заказе" (between runs)  -> E17 runner speed (s/op), 10 -> 8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Accelerate purchasers: reduce raw times (seconds) for the purchasing-related rows.
$ws.Range("D2").Value = 900
$ws.Range("D3").Value = 1800
$ws.Range("D4").Value = 900

# Accelerate process runner from 10 s to 8 s (D17 = E17*B18 recalculates automatically).
$ws.Range("E17").Value = 8

# Update the active selection to match the saved state of the workbook.
$ws.Range("F9").Select()

$wb.Save()
